$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text formatting so numeric-looking strings
# (e.g. "309.84", "0.3920", "11.30") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.037.06"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.814.24"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "309.84"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.5007"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("D8").Value = "0.3920"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "0.09878"
$ws.Range("E9").Value = "  +26.29%  "
$ws.Range("D10").Value = "1.104"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "40.88"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "6.415"
$ws.Range("E12").Value = "  +3.14%  "
$ws.Range("D13").Value = "20.55"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Value = "0.9996"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "1.810.24"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "7.279"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "0.00001137"
$ws.Range("E17").Value = "  +5.88%  "
$ws.Range("D18").Value = "92.43"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "0.06640"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "0.9985"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "17.18"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "5.927"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "28.103.93"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "158.63"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").Value = "2.021.85"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").Value = "2.415"
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("D30").Value = "127.11"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "5.573"
$ws.Range("E33").Value = "  +1.71%  "
$ws.Range("D34").Value = "3.593"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").Value = "0.06732"
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("D36").Value = "8.954"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").Value = "0.02331"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "0.2141"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "4.944"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "11.30"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "0.6189"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "1.173"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").Value = "0.9979"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "13.20"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "0.5911"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "3.690"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "124.31"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "0.06789"
$ws.Range("E51").Value = "  -0.29%  "
